$wb = $excel.ActiveWorkbook

# --- Sheet2: rename "GC-ECD-SF6" -> "GCECD" and trim it down to a single station (CGO) ---
$ws2 = $wb.Worksheets.Item("GC-ECD-SF6")
$ws2.Name = "GCECD"

# Row 7: only "General release date" label remains, drop the date value entirely
$ws2.Range("B7").Clear() | Out-Null

# Row 8: species header column now only has CGO (was MHD in col B); drop THD/RPB/SMO/CGO cols
$ws2.Range("B8").Value = "CGO"
$ws2.Range("C8:F8").ClearContents() | Out-Null

# Row 9: SF6 release date for CGO moves into column B; drop the other station columns
$ws2.Range("B9").Value = "2009-06-30 00:00"
$ws2.Range("C9:F9").ClearContents() | Out-Null

# Remove the now-unused columns C:F entirely
$ws2.Range("C:F").Delete() | Out-Null

# --- Sheet selection / active tab bookkeeping ---
$ws1 = $wb.Worksheets.Item("GCMD")
$ws1.Range("F21").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("E12").Select() | Out-Null
